$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values stored as literal text (e.g. "1.000", "0.7130").
# Force those cells to Text number format before writing so Excel does not
# reinterpret/round them as numeric values and strip significant characters.
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D33", "D35", "D36", "D39", "D40", "D41", "D42", "D43", "D44", "D47", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.323.12'
$ws.Range("E2").Value = '  +0.01%  '

$ws.Range("D3").Value = '1.876.61'
$ws.Range("E3").Value = '  +0.17%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '0.7130'
$ws.Range("E5").Value = '  +0.25%  '

$ws.Range("D6").Value = '242.85'
$ws.Range("E6").Value = '  +0.47%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = '0.08027'
$ws.Range("E8").Value = '  +3.42%  '

$ws.Range("D9").Value = '0.3158'
$ws.Range("E9").Value = '  +1.76%  '

$ws.Range("E10").Value = '  -0.04%  '

$ws.Range("D11").Value = '0.08240'
$ws.Range("E11").Value = '  -1.85%  '

$ws.Range("D12").Value = '1.893.74'
$ws.Range("E12").Value = '  +1.84%  '

$ws.Range("D13").Value = '5.255'

$ws.Range("D14").Value = '94.86'
$ws.Range("E14").Value = '  +4.18%  '

$ws.Range("D15").Value = '0.7129'

$ws.Range("D16").Value = '6.381'
$ws.Range("E16").Value = '  +5.31%  '

$ws.Range("D17").Value = '0.000008568'
$ws.Range("E17").Value = '  +4.77%  '

$ws.Range("D18").Value = '29.355.81'
$ws.Range("E18").Value = '  +0.12%  '

$ws.Range("D19").Value = '243.83'
$ws.Range("E19").Value = '  +1.89%  '

$ws.Range("D20").Value = '2.155.09'
$ws.Range("E20").Value = '  +1.67%  '

$ws.Range("D21").Value = '13.26'
$ws.Range("E21").Value = '  +0.53%  '

$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("E23").Value = '  +0.45%  '

$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  +0.04%  '

$ws.Range("D25").Value = '0.1563'
$ws.Range("E25").Value = '  -1.65%  '

$ws.Range("D26").Value = '9.049'
$ws.Range("E26").Value = '  +0.35%  '

$ws.Range("D27").Value = '162.60'
$ws.Range("E27").Value = '  -0.04%  '

$ws.Range("D28").Value = '18.54'
$ws.Range("E28").Value = '  +0.33%  '

$ws.Range("D29").Value = '1.505'
$ws.Range("E29").Value = '  -0.20%  '

$ws.Range("E30").Value = '  +0.52%  '

$ws.Range("D31").Value = '4.308'
$ws.Range("E31").Value = '  -0.14%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '0.05374'
$ws.Range("E32").Value = '  +1.61%  '

$ws.Range("B33").Value = 'Toncoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D33").Value = '1.179'
$ws.Range("E33").Value = '  -8.39%  '

$ws.Range("E34").Value = '  +0.30%  '

$ws.Range("D35").Value = '0.7652'
$ws.Range("E35").Value = '  +2.79%  '

$ws.Range("D36").Value = '1.180'
$ws.Range("E36").Value = '  +0.45%  '

$ws.Range("E37").Value = '  -0.46%  '

$ws.Range("E38").Value = '  +0.17%  '

$ws.Range("D39").Value = '1.253.41'
$ws.Range("E39").Value = '  +2.69%  '

$ws.Range("D40").Value = '2.753'
$ws.Range("E40").Value = '  +0.98%  '

$ws.Range("D41").Value = '6.481'
$ws.Range("E41").Value = '  -0.53%  '

$ws.Range("D42").Value = '0.9142'
$ws.Range("E42").Value = '  +3.19%  '

$ws.Range("D43").Value = '112.91'
$ws.Range("E43").Value = '  +3.22%  '

$ws.Range("D44").Value = '74.19'
$ws.Range("E44").Value = '  +2.55%  '

$ws.Range("E45").Value = '  +9.11%  '

$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("D47").Value = '2.045.05'
$ws.Range("E47").Value = '  +1.26%  '

$ws.Range("E48").Value = '  +0.50%  '

$ws.Range("D49").Value = '1.802'
$ws.Range("E49").Value = '  +0.32%  '

$ws.Range("D50").Value = '9.473'
$ws.Range("E50").Value = '  +1.31%  '

$ws.Range("D51").Value = '0.4356'
$ws.Range("E51").Value = '  +0.99%  '
